$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap F:V data between paired rows (home/away draws reordered) ---
function Swap-Rows($ws, $r1, $r2) {
    $tmp = $ws.Range("F$r1" + ":V$r1").Value2
    $ws.Range("F$r1" + ":V$r1").Value2 = $ws.Range("F$r2" + ":V$r2").Value2
    $ws.Range("F$r2" + ":V$r2").Value2 = $tmp
}

Swap-Rows $ws 19 20
Swap-Rows $ws 22 23
Swap-Rows $ws 32 33
Swap-Rows $ws 60 61

# --- Step 2: append new rows 72-79 ---
# Seed formatting for the new rows by copying the format of the last existing row (row 71)
$ws.Range("A71:V71").Copy()
$ws.Range("A72:V79").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 72
$ws.Range("A72").Value2 = 71
$ws.Range("B72").Value2 = "argentina"
$ws.Range("C72").Value2 = "copa-de-la-liga-profesional"
$ws.Range("D72").NumberFormat = "@"
$ws.Range("D72").Value2 = "2023"
$ws.Range("D72").Style = "Normal"
$ws.Range("E72").Value2 = 45192.875
$ws.Range("F72").Value2 = "Gimnasia L.P."
$ws.Range("G72").Value2 = 2
$ws.Range("H72").Value2 = "Rosario Central"
$ws.Range("I72").Value2 = 1
$ws.Range("J72").Value2 = 2.93
$ws.Range("K72").Value2 = "20/09/2023 01:12"
$ws.Range("L72").Value2 = 2.06
$ws.Range("M72").Value2 = "23/09/2023 20:46"
$ws.Range("N72").Value2 = 3.01
$ws.Range("O72").Value2 = "20/09/2023 01:12"
$ws.Range("P72").Value2 = 3.28
$ws.Range("Q72").Value2 = "23/09/2023 20:46"
$ws.Range("R72").Value2 = 2.71
$ws.Range("S72").Value2 = "20/09/2023 01:12"
$ws.Range("T72").Value2 = 4.09
$ws.Range("U72").Value2 = "23/09/2023 20:46"
$ws.Range("V72").Value2 = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/gimnasia-l-p-rosario-central/E1JwUDP3/"

# Row 73
$ws.Range("A73").Value2 = 72
$ws.Range("B73").Value2 = "argentina"
$ws.Range("C73").Value2 = "copa-de-la-liga-profesional"
$ws.Range("D73").NumberFormat = "@"
$ws.Range("D73").Value2 = "2023"
$ws.Range("D73").Style = "Normal"
$ws.Range("E73").Value2 = 45192.97916666666
$ws.Range("F73").Value2 = "Boca Juniors"
$ws.Range("G73").Value2 = 1
$ws.Range("H73").Value2 = "Lanus"
$ws.Range("I73").Value2 = 1
$ws.Range("J73").Value2 = 2.06
$ws.Range("K73").Value2 = "19/09/2023 23:12"
$ws.Range("L73").Value2 = 2.18
$ws.Range("M73").Value2 = "23/09/2023 23:29"
$ws.Range("N73").Value2 = 3.14
$ws.Range("O73").Value2 = "19/09/2023 23:12"
$ws.Range("P73").Value2 = 3.1
$ws.Range("Q73").Value2 = "23/09/2023 23:27"
$ws.Range("R73").Value2 = 3.86
$ws.Range("S73").Value2 = "19/09/2023 23:12"
$ws.Range("T73").Value2 = 3.96
$ws.Range("U73").Value2 = "23/09/2023 23:29"
$ws.Range("V73").Value2 = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/boca-juniors-lanus/vNeOuBOF/"

# Row 74
$ws.Range("A74").Value2 = 73
$ws.Range("B74").Value2 = "argentina"
$ws.Range("C74").Value2 = "copa-de-la-liga-profesional"
$ws.Range("D74").NumberFormat = "@"
$ws.Range("D74").Value2 = "2023"
$ws.Range("D74").Style = "Normal"
$ws.Range("E74").Value2 = 45193.08333333334
$ws.Range("F74").Value2 = "Central Cordoba"
$ws.Range("G74").Value2 = 2
$ws.Range("H74").Value2 = "Defensa y Justicia"
$ws.Range("I74").Value2 = 1
$ws.Range("J74").Value2 = 2.13
$ws.Range("K74").Value2 = "20/09/2023 20:12"
$ws.Range("L74").Value2 = 2.02
$ws.Range("M74").Value2 = "24/09/2023 01:52"
$ws.Range("N74").Value2 = 3.27
$ws.Range("O74").Value2 = "20/09/2023 20:12"
$ws.Range("P74").Value2 = 3.32
$ws.Range("Q74").Value2 = "24/09/2023 01:35"
$ws.Range("R74").Value2 = 3.5
$ws.Range("S74").Value2 = "20/09/2023 20:12"
$ws.Range("T74").Value2 = 4.23
$ws.Range("U74").Value2 = "24/09/2023 01:52"
$ws.Range("V74").Value2 = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/central-cordoba-santiago-del-estero-defensa-y-justicia/MooTvVwM/"

# Row 75
$ws.Range("A75").Value2 = 74
$ws.Range("B75").Value2 = "argentina"
$ws.Range("C75").Value2 = "copa-de-la-liga-profesional"
$ws.Range("D75").NumberFormat = "@"
$ws.Range("D75").Value2 = "2023"
$ws.Range("D75").Style = "Normal"
$ws.Range("E75").Value2 = 45193.08333333334
$ws.Range("F75").Value2 = "Huracan"
$ws.Range("G75").Value2 = 3
$ws.Range("H75").Value2 = "Velez Sarsfield"
$ws.Range("I75").Value2 = 0
$ws.Range("J75").Value2 = 2.27
$ws.Range("K75").Value2 = "20/09/2023 22:42"
$ws.Range("L75").Value2 = 2.49
$ws.Range("M75").Value2 = "24/09/2023 01:59"
$ws.Range("N75").Value2 = 2.92
$ws.Range("O75").Value2 = "20/09/2023 22:42"
$ws.Range("P75").Value2 = 2.71
$ws.Range("Q75").Value2 = "24/09/2023 01:59"
$ws.Range("R75").Value2 = 3.84
$ws.Range("S75").Value2 = "20/09/2023 22:42"
$ws.Range("T75").Value2 = 3.81
$ws.Range("U75").Value2 = "24/09/2023 01:59"
$ws.Range("V75").Value2 = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/huracan-velez-sarsfield/WxUYUgAc/"

# Row 76
$ws.Range("A76").Value2 = 75
$ws.Range("B76").Value2 = "argentina"
$ws.Range("C76").Value2 = "copa-de-la-liga-profesional"
$ws.Range("D76").NumberFormat = "@"
$ws.Range("D76").Value2 = "2023"
$ws.Range("D76").Style = "Normal"
$ws.Range("E76").Value2 = 45193.83333333334
$ws.Range("F76").Value2 = "Independiente"
$ws.Range("G76").Value2 = 0
$ws.Range("H76").Value2 = "Instituto"
$ws.Range("I76").Value2 = 0
$ws.Range("J76").Value2 = 2.14
$ws.Range("K76").Value2 = "21/09/2023 01:12"
$ws.Range("L76").Value2 = 2.35
$ws.Range("M76").Value2 = "24/09/2023 19:59"
$ws.Range("N76").Value2 = 3.02
$ws.Range("O76").Value2 = "21/09/2023 01:12"
$ws.Range("P76").Value2 = 2.91
$ws.Range("Q76").Value2 = "24/09/2023 19:58"
$ws.Range("R76").Value2 = 4.08
$ws.Range("S76").Value2 = "21/09/2023 01:12"
$ws.Range("T76").Value2 = 3.79
$ws.Range("U76").Value2 = "24/09/2023 19:59"
$ws.Range("V76").Value2 = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/independiente-instituto/hhIsTXvA/"

# Row 77
$ws.Range("A77").Value2 = 76
$ws.Range("B77").Value2 = "argentina"
$ws.Range("C77").Value2 = "copa-de-la-liga-profesional"
$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value2 = "2023"
$ws.Range("D77").Style = "Normal"
$ws.Range("E77").Value2 = 45193.94791666666
$ws.Range("F77").Value2 = "Newells Old Boys"
$ws.Range("G77").Value2 = 0
$ws.Range("H77").Value2 = "Estudiantes L.P."
$ws.Range("I77").Value2 = 1
$ws.Range("J77").Value2 = 2.35
$ws.Range("K77").Value2 = "21/09/2023 01:12"
$ws.Range("L77").Value2 = 2.39
$ws.Range("M77").Value2 = "24/09/2023 22:43"
$ws.Range("N77").Value2 = 2.94
$ws.Range("O77").Value2 = "21/09/2023 01:12"
$ws.Range("P77").Value2 = 2.8
$ws.Range("Q77").Value2 = "24/09/2023 22:43"
$ws.Range("R77").Value2 = 3.41
$ws.Range("S77").Value2 = "21/09/2023 01:12"
$ws.Range("T77").Value2 = 3.86
$ws.Range("U77").Value2 = "24/09/2023 22:44"
$ws.Range("V77").Value2 = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/newells-old-boys-estudiantes-l-p/W6eUJT1d/"

# Row 78
$ws.Range("A78").Value2 = 77
$ws.Range("B78").Value2 = "argentina"
$ws.Range("C78").Value2 = "copa-de-la-liga-profesional"
$ws.Range("D78").NumberFormat = "@"
$ws.Range("D78").Value2 = "2023"
$ws.Range("D78").Style = "Normal"
$ws.Range("E78").Value2 = 45194
$ws.Range("F78").Value2 = "Tigre"
$ws.Range("G78").Value2 = 0
$ws.Range("H78").Value2 = "San Lorenzo"
$ws.Range("I78").Value2 = 0
$ws.Range("J78").Value2 = 2.14
$ws.Range("K78").Value2 = "20/09/2023 22:42"
$ws.Range("L78").Value2 = 2.39
$ws.Range("M78").Value2 = "24/09/2023 23:53"
$ws.Range("N78").Value2 = 3.02
$ws.Range("O78").Value2 = "20/09/2023 22:42"
$ws.Range("P78").Value2 = 2.87
$ws.Range("Q78").Value2 = "24/09/2023 23:51"
$ws.Range("R78").Value2 = 3.8
$ws.Range("S78").Value2 = "20/09/2023 22:42"
$ws.Range("T78").Value2 = 3.75
$ws.Range("U78").Value2 = "24/09/2023 23:53"
$ws.Range("V78").Value2 = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/tigre-san-lorenzo/hGdQK9nj/"

# Row 79
$ws.Range("A79").Value2 = 78
$ws.Range("B79").Value2 = "argentina"
$ws.Range("C79").Value2 = "copa-de-la-liga-profesional"
$ws.Range("D79").NumberFormat = "@"
$ws.Range("D79").Value2 = "2023"
$ws.Range("D79").Style = "Normal"
$ws.Range("E79").Value2 = 45194.08333333334
$ws.Range("F79").Value2 = "Banfield"
$ws.Range("G79").Value2 = 1
$ws.Range("H79").Value2 = "River Plate"
$ws.Range("I79").Value2 = 1
$ws.Range("J79").Value2 = 3.57
$ws.Range("K79").Value2 = "22/09/2023 01:12"
$ws.Range("L79").Value2 = 3.89
$ws.Range("M79").Value2 = "25/09/2023 01:56"
$ws.Range("N79").Value2 = 3.29
$ws.Range("O79").Value2 = "22/09/2023 01:12"
$ws.Range("P79").Value2 = 3.36
$ws.Range("Q79").Value2 = "25/09/2023 01:56"
$ws.Range("R79").Value2 = 2.18
$ws.Range("S79").Value2 = "22/09/2023 01:12"
$ws.Range("T79").Value2 = 2.08
$ws.Range("U79").Value2 = "25/09/2023 01:56"
$ws.Range("V79").Value2 = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/banfield-river-plate/MX5fQVOS/"
